$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "columns": renumber the 4.3.2/4.4/4.4.1/4.5 questions to
# 4.4/4.5/4.5.1/4.6 (an extra question was inserted earlier in the form,
# shifting the numbering of pregMonths, lmpKnown, lmpDate and
# expectedDeliveryDate down by one full step) and fix the eddDate
# calculation-related formatting.
# ---------------------------------------------------------------------------
$wsColumns = $wb.Worksheets.Item("columns")

# pregMonths (row 12): 4.3.2. -> 4.4.
$wsColumns.Range("H12").Value = "4.4. Approximately how many completed months has the woman been pregnant?"
$wsColumns.Range("I12").Value = "4.4. Número aproximado de meses da gravidez?"
$wsColumns.Range("J12").Value = "4.4. Nombre approximatif de mois de grossesse?"

# lmpKnown (row 13): 4.4. -> 4.5.
$wsColumns.Range("H13").Value = "4.5. Is the date of start of the last menstrual period known?"
$wsColumns.Range("I13").Value = "4.5. Conhece a data em que iniciou a sua última menstruação?"
$wsColumns.Range("J13").Value = "4.5. Connaissez-vous la date du début de vos dernières règles?"

# lmpDate (row 14): 4.4.1. -> 4.5.1.
$wsColumns.Range("H14").Value = "4.5.1. Date of last menstrual period"
$wsColumns.Range("I14").Value = "4.5.1. Data em que iniciou a última menstruação"
$wsColumns.Range("J14").Value = "4.5.1. Date du début des dernières règles"

# expectedDeliveryDate (row 15): 4.5. -> 4.6.
$wsColumns.Range("H15").Value = "4.6. Calculated Expected Delivery Date"
$wsColumns.Range("I15").Value = "4.6. Data provável do parto pré-calculada"
$wsColumns.Range("J15").Value = "4.6. Date probable d'accouchement pré-calculée"

# display_condition column (O) for hasPrenatalRecord/eddDate/eddType rows is
# formula-looking text ("${status}='PREGNANT' and ${eddKnown}='TRUE'") -
# force it to be stored/displayed as plain text, matching the rest of the
# column.
$wsColumns.Range("O9:O11").NumberFormat = "@"

# Restore the view: scrolled right so column I is the first visible column,
# with the active cell on N15.
$wsColumns.Activate()
$excel.ActiveWindow.ScrollColumn = 9
$wsColumns.Range("N15").Select()

# ---------------------------------------------------------------------------
# Sheet "options": no data changed, just the last selected cell.
# ---------------------------------------------------------------------------
$wsOptions = $wb.Worksheets.Item("options")
$wsOptions.Activate()
$wsOptions.Range("B2").Select()

# Re-activate the "columns" sheet so it stays the selected tab, matching the
# original workbook (tabSelected="1" on the columns sheet).
$wsColumns.Activate()
